# Apply the data changes to row 2 of Sheet1.
# These cells store plain numeric-looking text (inlineStr in the OOXML,
# i.e. Excel's "number stored as text"), so a direct `.Value = "..."`
# assignment would be auto-coerced to a real Number by Excel and would
# also pick up a stray style (quote-prefix / text-format) on the cell.
# To faithfully reproduce "text that happens to look like a number"
# without touching cell formatting, stage each new value as text in a
# scratch cell (using a leading apostrophe to force text), copy it, and
# PasteSpecial just the values onto the destination cell. The scratch
# cell is cleared afterwards so it leaves no trace in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$xlPasteValues = -4163

function Set-TextValue($cellRef, $text) {
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial($xlPasteValues)
}

Set-TextValue "B2" "1.0"
Set-TextValue "C2" "8761"
Set-TextValue "E2" "0.06"
Set-TextValue "F2" "0.661"
Set-TextValue "P2" "3.366"
Set-TextValue "Q2" "50346.16"

$scratch.Clear()
